$wb = $excel.ActiveWorkbook

# Rename sheets
$wb.Worksheets.Item(1).Name = "GNG_TO-16504778668283727"
$wb.Worksheets.Item(2).Name = "NB_TO-16504778686756663"
$wb.Worksheets.Item(3).Name = "RS_TO-16504778686816652"
$wb.Worksheets.Item(4).Name = "TOL_TO-16504778687406635"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16504778688016996"

# Sheet 1 (GNG_TO)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16504778667873728.csv"
$ws1.Range("B3").Value = "GNG_stims-1650477866811456.csv"
$ws1.Range("B4").Value = "go_stims-16504778668133724.csv"
$ws1.Range("B5").Value = "GNG_stims-16504778668273969.csv"

# Sheet 2 (NB_TO)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "TB-1650477868305665.csv"
$ws2.Range("B3").Value = "ZB-match_7-16504778672984045.csv"
$ws2.Range("B4").Value = "ZB-match_0-16504778672113698.csv"
$ws2.Range("B5").Value = "OB-16504778673300564.csv"
$ws2.Range("B6").Value = "TB-16504778686586974.csv"
$ws2.Range("B7").Value = "ZB-match_5-16504778669234076.csv"
$ws2.Range("B8").Value = "OB-16504778675333512.csv"
$ws2.Range("B9").Value = "TB-1650477868121696.csv"
$ws2.Range("B10").Value = "OB-16504778680686617.csv"

# Sheet 3 (RS_TO) - no data changes

# Sheet 4 (TOL_TO)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16504778687076976.csv"
$ws4.Range("B3").Value = "ZM_stims-1650477868684666.csv"
$ws4.Range("B4").Value = "MM_stims-16504778687236981.csv"
$ws4.Range("B5").Value = "ZM_stims-16504778687086635.csv"
$ws4.Range("B6").Value = "MM_stims-16504778687396648.csv"
$ws4.Range("B7").Value = "ZM_stims-16504778687246654.csv"

# Sheet 5 (vSAT_TO)
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-16504778687856622.csv"
$ws5.Range("B3").Value = "SAT_stims-1650477868743691.csv"
$ws5.Range("B4").Value = "SAT_stims-16504778687556887.csv"
$ws5.Range("B5").Value = "vSAT_stims-16504778687706985.csv"
